$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Update SMD designator lists / quantities for existing rows ---

# Row 7: capacitors ">= 25V" now also cover C14,C15
$ws.Range("A7").Value = "C14,C15,C19,C20,C21,C22,C23,C24,C25,C26,C27,C28"

# Row 10: diodes D1,D2 -> D1,D2,D4 (qty 2 -> 3)
$ws.Range("A10").Value = "D1,D2,D4"
$ws.Range("E10").Value = 3

# Row 14: resistors R15,R30,R32,R39 -> +R44 (qty 4 -> 5)
$ws.Range("A14").Value = "R15,R30,R32,R39,R44"
$ws.Range("E14").Value = 5

# --- 2) Insert new SMD part row for U6 (single op-amp), right above the
#        THT-parts section header (old row 30) ---
$ws.Range("A30").EntireRow.Insert()
$ws.Range("A30").Value = "U6"
$ws.Range("B30").Value = "Single op-amp"
$ws.Range("C30").Value = "SOT-23-5"
$ws.Range("E30").Value = 1
$ws.Range("G30").Value = "Texas Instruments"
$ws.Range("H30").Value = "LM321MFX/NOPB"
$ws.Range("I30").Value = "926-LM321MFX/NOPB"

# --- 3) Insert new THT part row for D6 (red LED, input indicator), as the
#        first data row of the THT-parts section (new row 32, right below
#        the THT header which just shifted from row 30 to row 31) ---
$ws.Range("A32").EntireRow.Insert()
# reset formatting inherited from the bold section-header row above so the
# new row matches the plain "normal" style used by the other data rows
$ws.Range("A32:I32").Font.Bold = $false
$ws.Range("A32:I32").Interior.ColorIndex = 0
$ws.Range("A32").Value = "D6"
$ws.Range("B32").Value = "LED, red, 5mm"
$ws.Range("C32").Value = "T-1"
$ws.Range("E32").Value = 1
$ws.Range("G32").Value = "Kingbright"
$ws.Range("H32").Value = "WP483IDT"
$ws.Range("I32").Value = "604-WP483IDT"

# --- 4) Fix up the two section-total formulas so their ranges cover the
#        newly inserted rows (Excel's auto-grow only extends a SUM() range
#        when the insertion point is strictly inside it, not right at the
#        boundary row, so these need to be set explicitly) ---
$ws.Range("A3").Formula = '=_xlfn.CONCAT("SMD Parts: ",SUM(E4:E30)," pcs")'
$ws.Range("A31").Formula = '=_xlfn.CONCAT("THT Parts: ",SUM(E32:E40)," pcs")'
